# "Big update for new student year"
# Replace outgoing faculty member (陈锶奇 / Siqi Chen) in row 4 with the
# incoming one (马亿 / Yi Ma), and tweak Tianpei Yang's title in row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: swap in the new faculty member's info ---------------------
$ws.Range("A4").Value = "马亿"
$ws.Range("B4").Value = "Yi Ma"
$ws.Range("C4").Value = "Associate Professor,  School of Comuputer and Information Technology (School of Big Data), Shanxi University"
$ws.Range("D4").Value = "Reinforcement Learning, Embodied AI, RL for Application"
$ws.Range("E4").Value = "https://mayi1996.top/"

# New title text is long, so wrap it and grow the row to fit.
# (Applied before the hyperlink rebuild below so this wrap-text style lands
# at the next free style slot, ahead of the style churn Hyperlinks.Add()
# causes.)
$ws.Range("C4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 33

# --- Row 5: tweak Tianpei Yang's position text --------------------------
$ws.Range("C5").Value = "Postdoc at University of Alberta; "

# --- Fix up the personal-homepage hyperlink for row 4 --------------------
# (Range.Hyperlinks.Delete() removes every hyperlink on the sheet in this
# engine, so capture the addresses we want to keep first and re-add them,
# then restore their original "hyperlink" cell style.)
$addrs = @()
foreach ($h in $ws.Hyperlinks) {
    $addrs += $h.Address
}
$ws.Range("E4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), $addrs[0])
$ws.Hyperlinks.Add($ws.Range("E3"), $addrs[1])
$ws.Hyperlinks.Add($ws.Range("E4"), "https://mayi1996.top/")
$ws.Range("E2").Style = "超链接"
$ws.Range("E3").Style = "超链接"
$ws.Range("E4").Style = "超链接"

# --- Misc view state ------------------------------------------------------
$ws.Range("C5").Select()
